$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date text (slash -> hyphen)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Rows whose new date text is ambiguous to Excel's smart date parser
# (day-of-month <= 12, so "DD-MM-YYYY" could also read as "MM-DD-YYYY") need
# to be forced to stay literal text, otherwise Excel auto-converts them to
# a real date serial number.
$ambiguousRows = @(4, 5, 6, 7, 13, 14, 15, 16)

foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    if ($ambiguousRows -contains $r) {
        $cell.NumberFormat = "@"
        $cell.Value = $dates[$r]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $dates[$r]
    }
}

# Rows where D and E go 0 -> 1, and H goes 1 -> 0
$flipRows = @(4, 5, 6, 10, 11, 13, 14)

foreach ($r in $flipRows) {
    $ws.Cells.Item($r, 4).Value = 1   # D
    $ws.Cells.Item($r, 5).Value = 1   # E
    $ws.Cells.Item($r, 8).Value = 0   # H
}
